$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 65; A = "CK_CU_BOUNDARY_En"; B = "CU_Name"; C = "Consistency (C1)"; D = "2024-12-02 23:10:16"; E = 0.91; F = 1; G = "OnakD" },
    @{ Row = 66; A = "CK_CU_BOUNDARY_En"; B = "Shape_Length, Shape_Area"; C = "Accuracy (A1)"; D = "2024-12-02 23:10:16"; E = "no threshold"; F = 0.993670886075949; G = "OnakD" },
    @{ Row = 67; A = "Conservation_Unit_Data_20220902"; B = "ACT_ID, ANALYSIS_YR, STREAM_ID, SPL_ID, NATURAL_ADULT_SPAWNERS, NATURAL_JACK_SPAWNERS, JACK_BROODSTOCK_REMOVALS, TOTAL_BROODSTOCK_REMOVALS, OTHER_REMOVALS, TOTAL_RETURN_TO_RIVER, UNSPECIFIED_RETURN, NO_INSPECTIONS_USED, MAX_ESTIMATE, EFFECTIVE_FEMALES, WEIGHTED_PCT_SPAWN, OTHER_ADULT_REMOVALS, OTHER_JACK_REMOVALS, TOT_ADULT_RET_RIVER, POP_ID, SBJ_ID"; C = "Accuracy (A1)"; D = "2024-12-03 09:00:49"; E = "no threshold"; F = 1; G = "onakd" },
    @{ Row = 68; A = "Conservation_Unit_Data_20220902"; B = "All columns"; C = "Accuracy (A3)"; D = "2024-12-03 09:02:20"; E = "no threshold"; F = 1; G = "onakd" },
    @{ Row = 69; A = "Conservation_Unit_Data_20220902"; B = "All columns"; C = "Completeness (P)"; D = "2024-12-03 09:02:40"; E = 0.75; F = 0.924034635876363; G = "onakd" },
    @{ Row = 70; A = "Conservation_Unit_Data_20220902"; B = "ACT_ID, ANALYSIS_YR, STREAM_ID, SPL_ID, NATURAL_ADULT_SPAWNERS, NATURAL_JACK_SPAWNERS, JACK_BROODSTOCK_REMOVALS, TOTAL_BROODSTOCK_REMOVALS, OTHER_REMOVALS, TOTAL_RETURN_TO_RIVER, UNSPECIFIED_RETURN, NO_INSPECTIONS_USED, MAX_ESTIMATE, EFFECTIVE_FEMALES, WEIGHTED_PCT_SPAWN, OTHER_ADULT_REMOVALS, OTHER_JACK_REMOVALS, TOT_ADULT_RET_RIVER, TOT_JACK_RET_RIVER, JUV_PRES_TYP, POP_ID, SBJ_ID"; C = "Accuracy (A1)"; D = "2024-12-03 09:07:44"; E = "no threshold"; F = 0.9965248818459828; G = "onakd" }
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
}
